# Reorders/edits the Set Pair / Word Pair rows and recolors the Word Pair
# column fonts to match the refreshed 'colored_word_pairs' dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-BgrColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Chart.js-style palette used for the Word Pair font colors
$color_36A2EB = Get-BgrColor 0x36 0xA2 0xEB
$color_4BC0C0 = Get-BgrColor 0x4B 0xC0 0xC0
$color_9966FF = Get-BgrColor 0x99 0x66 0xFF
$color_FF6384 = Get-BgrColor 0xFF 0x63 0x84
$color_FF9F40 = Get-BgrColor 0xFF 0x9F 0x40
$color_FFCE56 = Get-BgrColor 0xFF 0xCE 0x56

$ws.Cells.Item(2, 1).Value = 'Bajta - Britto_2017'
$ws.Cells.Item(2, 2).Value = 'reliability - robustness level - maintainability level - reliability level - other project activities - type of project - not considered experience - in-house experience - experience level - tool experience level - oo experience level - domain experience level - software development experience - programming language experience level - effort hours - time efficiency level - maintenance - design - architecture - risk - risk level - considered structure - infrastructure - structure - value of number of sites - node count - considered experience - hardware - technical factors - platform support level - other requirement - requirements clarity level - requirements novelty level - processing requirements - security - security level - testing - testability level - not considered structure - maintainability - reusability level - installability level - portability level - size report - team size - availability - availability level - performance - quality level - process efficiency level - number of team members - team capability - work team level'
$ws.Cells.Item(2, 2).Font.Bold = $false
$ws.Cells.Item(2, 2).Font.Color = $color_9966FF

$ws.Cells.Item(3, 1).Value = 'Bajta - Britto_2016'
$ws.Cells.Item(3, 2).Value = 'temporal distance - geographic distance - socio-cultural distance - geographical distance'
$ws.Cells.Item(3, 2).Font.Bold = $true
$ws.Cells.Item(3, 2).Font.Color = $color_FFCE56

$ws.Cells.Item(4, 1).Value = 'Bajta - Dashti'
$ws.Cells.Item(4, 2).Value = 'expert judgment - machine learning - artificial neural networks - fuzzy similar - fuzzy logic - actual cost value - constructive cost model - estimated cost value'
$ws.Cells.Item(4, 2).Font.Bold = $true
$ws.Cells.Item(4, 2).Font.Color = $color_9966FF

$ws.Cells.Item(5, 1).Value = 'Bajta - Mendes'
$ws.Cells.Item(5, 2).Value = 'other requirement - specific - implementation - functionality - individual - nonspecific - size report - length - performance'
$ws.Cells.Item(5, 2).Font.Bold = $false
$ws.Cells.Item(5, 2).Font.Color = $color_FF9F40

$ws.Cells.Item(6, 1).Value = 'Bajta - Usman'
$ws.Cells.Item(6, 2).Value = 'reliability - maintainability - maintenance - availability - other project activities - all estimated activities - expert judgment - expert judgement - distant onshore - distributed: near offshore - distributed: close onshore - distributed: far offshore - distributed: distant onshore - preliminary planning - bidding planning level - release planning level - sprint planning level - daily planning level - planning poker - not considered experience - considered team''s prior experience - not considered team''s prior experience - effort hours - actual effort value - other effort predictors - ideal hours - not used effort predictors - telecommunication - communications industry - analysis - performance - near offshore - design - close onshore - considered structure - other type - agile - scrum - customized scrum - healthcare - health - considered experience - other requirement - other non functional requirements - not considered non functional requirements - other unit - far offshore - other planning approaches - execution - implementation - detail planning - system investigation - case-based reasoning - user case points - security - individual - single type - testing - actual cost value - group-based estimation - other estimation technique - number of team members - no. of team members - not considered team''s skill level - considered team''s skill level - considered distributed teams - not considered distributed teams - not applicable distributed teams - estimated cost value - estimate value(s) - finance - financial - statistics analysis'
$ws.Cells.Item(6, 2).Font.Bold = $true
$ws.Cells.Item(6, 2).Font.Color = $color_4BC0C0

$ws.Cells.Item(7, 1).Value = 'Britto_2017 - Bajta'
$ws.Cells.Item(7, 2).Value = 'risk level - risk - requirements clarity level - other requirement - in-house experience - not considered experience - considered experience - type of project - other project activities - robustness level - reliability - availability level - availability - technical factors - hardware - experience level - infrastructure - considered structure - testability level - testing - tool experience level - oo experience level - reusability level - maintainability - maintainability level - maintenance - node count - value of number of sites - requirements novelty level - time efficiency level - effort hours - performance - structure - not considered structure - team capability - number of team members - security level - security - domain experience level - reliability level - platform support level - architecture - design - software development experience - quality level - installability level - process efficiency level - team size - size report - programming language experience level - work team level - portability level - processing requirements'
$ws.Cells.Item(7, 2).Font.Bold = $false
$ws.Cells.Item(7, 2).Font.Color = $color_4BC0C0

$ws.Cells.Item(8, 1).Value = 'Britto_2017 - Mendes'
$ws.Cells.Item(8, 2).Value = 'cyclomatic complexity - complexity - web page allocation - web application - new complexity - page complexity - new media count - media - media duration - layout complexity - component complexity - metrics’ program - problem-oriented metric - media count - web objects - web software application - web hypermedia application - data usage complexity - total complexity - difficulty level - interface complexity - cohesion complexity - model collection complexity - model association complexity - control flow complexity - motivation level - motivation - spi program - program/sript - media allocation - lines of code - input complexity - output complexity - class complexity - data flow complexity - adaptation complexity'
$ws.Cells.Item(8, 2).Font.Bold = $false
$ws.Cells.Item(8, 2).Font.Color = $color_4BC0C0

$ws.Cells.Item(9, 1).Value = 'Britto_2017 - Usman'
$ws.Cells.Item(9, 2).Value = 'object-oriented heuristic function points - function points - requirements clarity level - other non functional requirements - not considered non functional requirements - type of project - other project domain - co-located project setting - use case count - use case points method - user case points - use case estimation entity - robustness level - reliability - availability level - availability - experience level - considered team''s prior experience - not considered team''s prior experience - productivity level - other effort predictors - communication level - considered customer communication - not considered customer communication - type of product - other type - single type - requirements volatility level - testability level - testing - reusability level - maintainability - maintainability level - maintenance - requirements novelty level - time efficiency level - performance - actual effort value - international function point users group - team capability - no. of team members - not considered team''s skill level - considered team''s skill level - considered distributed teams - not considered distributed teams - not applicable distributed teams - security level - security - publishing unit count - other unit - reliability level - object-oriented function points - point type - entity count - number of entities estimated - architecture - design - quality level - installability level - implementation - process efficiency level - team size - work team level - portability level - processing requirements'
$ws.Cells.Item(9, 2).Font.Bold = $false
$ws.Cells.Item(9, 2).Font.Color = $color_4BC0C0

$ws.Cells.Item(10, 1).Value = 'Britto_2016 - Bajta'
$ws.Cells.Item(10, 2).Value = 'temporal distance - geographical distance - geographic distance - socio-cultural distance'
$ws.Cells.Item(10, 2).Font.Bold = $true
$ws.Cells.Item(10, 2).Font.Color = $color_36A2EB

$ws.Cells.Item(11, 1).Value = 'Britto_2016 - Dashti'
$ws.Cells.Item(11, 2).Value = 'legal entity'
$ws.Cells.Item(11, 2).Font.Bold = $true
$ws.Cells.Item(11, 2).Font.Color = $color_FF9F40

$ws.Cells.Item(12, 1).Value = 'Britto_2016 - Mendes'
$ws.Cells.Item(12, 2).Value = 'late estimation stage - early - early estimation stage - distributed - direct - early & late estimation stage'
$ws.Cells.Item(12, 2).Font.Bold = $false
$ws.Cells.Item(12, 2).Font.Color = $color_FF9F40

$ws.Cells.Item(13, 1).Value = 'Britto_2016 - Usman'
$ws.Cells.Item(13, 2).Value = 'late estimation stage - other estimation technique - other estimation entity - estimate value(s) - early estimation stage - estimator & provider - use case estimation entity - task estimation entity - estimator - semi-distributed - distribution type - distributed: near offshore - distributed: far offshore - distributed - considered distributed teams - not considered distributed teams - not applicable distributed teams - distributed: close onshore - distributed: distant onshore - early & late estimation stage'
$ws.Cells.Item(13, 2).Font.Bold = $false
$ws.Cells.Item(13, 2).Font.Color = $color_FF9F40

$ws.Cells.Item(14, 1).Value = 'Dashti - Bajta'
$ws.Cells.Item(14, 2).Value = 'fuzzy logic - fuzzy similar - expert judgment - constructive cost model - actual cost value - estimated cost value - artificial neural networks - machine learning'
$ws.Cells.Item(14, 2).Font.Bold = $true
$ws.Cells.Item(14, 2).Font.Color = $color_9966FF

$ws.Cells.Item(15, 1).Value = 'Dashti - Britto_2016'
$ws.Cells.Item(15, 2).Value = 'legal entity'
$ws.Cells.Item(15, 2).Font.Bold = $true
$ws.Cells.Item(15, 2).Font.Color = $color_FF9F40

$ws.Cells.Item(16, 1).Value = 'Dashti - Usman'
$ws.Cells.Item(16, 2).Value = 'analogy-based - analogy - expert judgment - expert judgement'
$ws.Cells.Item(16, 2).Font.Bold = $false
$ws.Cells.Item(16, 2).Font.Color = $color_FF6384

$ws.Cells.Item(17, 1).Value = 'Mendes - Bajta'
$ws.Cells.Item(17, 2).Value = 'specific - other requirement - individual - length - size report - functionality - implementation - performance - nonspecific'
$ws.Cells.Item(17, 2).Font.Bold = $false
$ws.Cells.Item(17, 2).Font.Color = $color_4BC0C0

$ws.Cells.Item(18, 1).Value = 'Mendes - Britto_2017'
$ws.Cells.Item(18, 2).Value = 'web application - web page allocation - web objects - motivation - motivation level - web software application - web hypermedia application - media - new media count - media duration - media count - media allocation - program/sript - spi program - lines of code - problem-oriented metric - metrics’ program - complexity - cyclomatic complexity - new complexity - page complexity - layout complexity - component complexity - data usage complexity - total complexity - difficulty level - interface complexity - cohesion complexity - model collection complexity - model association complexity - control flow complexity - input complexity - output complexity - class complexity - data flow complexity - adaptation complexity'
$ws.Cells.Item(18, 2).Font.Bold = $false
$ws.Cells.Item(18, 2).Font.Color = $color_4BC0C0

$ws.Cells.Item(19, 1).Value = 'Mendes - Britto_2016'
$ws.Cells.Item(19, 2).Value = 'direct - distributed - early - late estimation stage - early estimation stage - early & late estimation stage'
$ws.Cells.Item(19, 2).Font.Bold = $false
$ws.Cells.Item(19, 2).Font.Color = $color_FF9F40

$ws.Cells.Item(20, 1).Value = 'Mendes - Usman'
$ws.Cells.Item(20, 2).Value = 'specific - other type - single type - solution-oriented metric - considered without any metric - problem-oriented metric - functionality - performance - other non functional requirements - not considered non functional requirements - implementation - nonspecific'
$ws.Cells.Item(20, 2).Font.Bold = $false
$ws.Cells.Item(20, 2).Font.Color = $color_4BC0C0

$ws.Cells.Item(21, 1).Value = 'Usman - Bajta'
$ws.Cells.Item(21, 2).Value = 'financial - finance - other type - considered structure - other requirement - other estimation technique - group-based estimation - considered team''s prior experience - not considered experience - considered experience - number of team members - security - communications industry - telecommunication - bidding planning level - preliminary planning - performance - analysis - implementation - release planning level - no. of team members - actual effort value - effort hours - actual cost value - estimated cost value - other non functional requirements - not considered team''s skill level - single type - individual - sprint planning level - other planning approaches - detail planning - considered team''s skill level - considered distributed teams - other effort predictors - scrum - agile - expert judgement - expert judgment - not considered distributed teams - daily planning level - system investigation - statistics analysis - design - not applicable distributed teams - not considered non functional requirements - customized scrum - distributed: near offshore - distant onshore - near offshore - close onshore - far offshore - execution - user case points - case-based reasoning - maintainability - reliability - maintenance - availability - ideal hours - distributed: close onshore - planning poker - other unit - not considered team''s prior experience - estimate value(s) - distributed: far offshore - distributed: distant onshore - testing - all estimated activities - other project activities - not used effort predictors - health - healthcare'
$ws.Cells.Item(21, 2).Font.Bold = $true
$ws.Cells.Item(21, 2).Font.Color = $color_4BC0C0

$ws.Cells.Item(22, 1).Value = 'Usman - Britto_2017'
$ws.Cells.Item(22, 2).Value = 'other type - type of product - considered team''s prior experience - experience level - team capability - work team level - security - security level - function points - object-oriented heuristic function points - international function point users group - object-oriented function points - performance - time efficiency level - quality level - process efficiency level - no. of team members - team size - actual effort value - other non functional requirements - requirements clarity level - requirements volatility level - requirements novelty level - processing requirements - not considered team''s skill level - single type - considered team''s skill level - considered distributed teams - considered customer communication - communication level - number of entities estimated - entity count - other effort predictors - productivity level - point type - not considered distributed teams - design - architecture - not applicable distributed teams - not considered non functional requirements - use case points method - use case count - implementation - installability level - user case points - maintainability - reusability level - maintainability level - reliability level - portability level - not considered customer communication - reliability - robustness level - maintenance - use case estimation entity - other unit - publishing unit count - not considered team''s prior experience - other project domain - type of project - testing - testability level - availability - availability level - co-located project setting'
$ws.Cells.Item(22, 2).Font.Bold = $false
$ws.Cells.Item(22, 2).Font.Color = $color_36A2EB

$ws.Cells.Item(23, 1).Value = 'Usman - Britto_2016'
$ws.Cells.Item(23, 2).Value = 'other estimation technique - late estimation stage - early estimation stage - estimator & provider - estimator - early & late estimation stage - distribution type - semi-distributed - distributed - other estimation entity - considered distributed teams - not considered distributed teams - not applicable distributed teams - distributed: near offshore - distributed: close onshore - use case estimation entity - estimate value(s) - distributed: far offshore - task estimation entity - distributed: distant onshore'
$ws.Cells.Item(23, 2).Font.Bold = $false
$ws.Cells.Item(23, 2).Font.Color = $color_FF9F40

$ws.Cells.Item(24, 1).Value = 'Usman - Dashti'
$ws.Cells.Item(24, 2).Value = 'analogy - analogy-based - expert judgement - expert judgment'
$ws.Cells.Item(24, 2).Font.Bold = $false
$ws.Cells.Item(24, 2).Font.Color = $color_4BC0C0

$ws.Cells.Item(25, 1).Value = 'Usman - Mendes'
$ws.Cells.Item(25, 2).Value = 'other type - specific - nonspecific - performance - functionality - other non functional requirements - single type - not considered non functional requirements - implementation - considered without any metric - solution-oriented metric - problem-oriented metric'
$ws.Cells.Item(25, 2).Font.Bold = $false
$ws.Cells.Item(25, 2).Font.Color = $color_9966FF
